$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.160.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.481.66"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.65"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.01"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.480.11"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.14%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.18"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -7.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.381"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -6.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.067.76"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.64"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.462.78"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.989.96"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.65"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -10.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.76"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.25%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "387.41"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -8.17%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.40"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.59%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.612.22"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.12%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.37"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.89%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -8.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.22"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -9.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.492.44"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.01"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "171.94"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -9.75%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -8.99%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -8.36%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -9.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0772"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.813"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.11"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -7.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.34"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -12.80%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.77"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.55%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.62"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -8.22%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.13"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.64"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.222.90"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.70%  "
